# started reformulating - not working yet (vardeviceflow instead of dispatchin/out)
#
# Adds a new "devicemodel" lookup sheet (between "device" and "hub") that maps a
# device model id to its "in"/"out" media types, and tags every existing device
# row on the "device" sheet with a model id in column P ("default" unless a row
# already specifies one, e.g. the compressor_el row).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "devicemodel" worksheet right before "hub" -----------
$hub = $wb.Worksheets.Item("hub")
$modelSheet = $wb.Worksheets.Add($hub)
$modelSheet.Name = "devicemodel"

# Fill it in roughly the order the strings were first typed so shared-string
# ids land the same way they would interactively.
$modelSheet.Range("A3").Value = "source_el"
$modelSheet.Range("A4").Value = "sink_el"
$modelSheet.Range("B1").Value = "in"
$modelSheet.Range("C1").Value = "out"
$modelSheet.Range("A5").Value = "gasturbine"
$modelSheet.Range("B2").Value = "el,gas"

$modelSheet.Range("A1").Value = "id"
$modelSheet.Range("A2").Value = "compressor_el"
$modelSheet.Range("C2").Value = "gas"
$modelSheet.Range("C3").Value = "el"
$modelSheet.Range("B4").Value = "el"
$modelSheet.Range("B5").Value = "gas"
$modelSheet.Range("C5").Value = "el"

$modelSheet.Range("C3").Select()

# --- 2. Tag every device row with its model id in column P ------------------
$deviceSheet = $wb.Worksheets.Item("device")

for ($row = 2; $row -le 15; $row++) {
    $cell = $deviceSheet.Cells.Item($row, 16)   # column P
    if ($cell.Text -eq "") {
        $cell.Value = "default"
    }
}

$deviceSheet.Range("P2").Select()
